# Auto-generated Excel COM-interop script
# Applies market-price data refresh values to Spriggan_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 657.4286
$ws.Range("I11").Value = 657.4286
$ws.Range("K11").Value = 657.4286
$ws.Range("M11").Value = -517.4286

$ws.Range("H40").Value = 3835842.8
$ws.Range("I40").Value = 3311.0952
$ws.Range("J40").Value = 13896239
$ws.Range("K40").Value = 3311.0952
$ws.Range("L40").Value = 13896239
$ws.Range("M40").Value = -3136.0952
$ws.Range("N40").Value = -13896589

$ws.Range("H43").Value = 5546.7334
$ws.Range("I43").Value = 2904.5454
$ws.Range("J43").Value = 12812.75
$ws.Range("K43").Value = 2904.5454
$ws.Range("L43").Value = 12812.75
$ws.Range("M43").Value = -2835.5454
$ws.Range("N43").Value = -12950.75

$ws.Range("H86").Value = 6585.125
$ws.Range("I86").Value = 7870.273
$ws.Range("J86").Value = 3757.8
$ws.Range("K86").Value = 7870.273
$ws.Range("L86").Value = 3757.8
$ws.Range("M86").Value = -6747.273
$ws.Range("N86").Value = -6003.8

$ws.Range("H89").Value = 6585.125
$ws.Range("I89").Value = 7870.273
$ws.Range("J89").Value = 3757.8
$ws.Range("K89").Value = 39351.365
$ws.Range("L89").Value = 18789
$ws.Range("M89").Value = -33735.365
$ws.Range("N89").Value = -30021

$ws.Range("H92").Value = 592.2727
$ws.Range("I92").Value = 551.4
$ws.Range("K92").Value = 551.4
$ws.Range("M92").Value = 696.6

$ws.Range("H94").Value = 2608.8
$ws.Range("I94").Value = 2586.25
$ws.Range("K94").Value = 2586.25
$ws.Range("M94").Value = -2135.25

$ws.Range("H137").Value = 2344.8572
$ws.Range("I137").Value = 2234.9
$ws.Range("J137").Value = 2619.75
$ws.Range("K137").Value = 6704.700000000001
$ws.Range("L137").Value = 7859.25
$ws.Range("M137").Value = -4154.700000000001
$ws.Range("N137").Value = -12959.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4039.4893
$ws.Range("I32").Value = 2175.1667
$ws.Range("K32").Value = 2175.1667
$ws.Range("M32").Value = -1888.1667

$ws.Range("H44").Value = 65000
$ws.Range("J44").Value = 65000
$ws.Range("L44").Value = 65000
$ws.Range("N44").Value = -65976

$ws.Range("H87").Value = 35000
$ws.Range("J87").Value = 35000
$ws.Range("L87").Value = 35000
$ws.Range("N87").Value = -37496

$ws.Range("H90").Value = 35000
$ws.Range("J90").Value = 35000
$ws.Range("L90").Value = 105000
$ws.Range("N90").Value = -117480

$ws.Range("H102").Value = 8339418.5
$ws.Range("I102").Value = 8339418.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 8339418.5
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -8337796.5
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 17242138
$ws.Range("I80").Value = 521
$ws.Range("K80").Value = 521
$ws.Range("M80").Value = 477

$ws.Range("H83").Value = 17242138
$ws.Range("I83").Value = 521
$ws.Range("K83").Value = 2605
$ws.Range("M83").Value = 2387

$ws.Range("H94").Value = 1246.3684
$ws.Range("I94").Value = 1320.1333
$ws.Range("J94").Value = 969.75
$ws.Range("K94").Value = 1320.1333
$ws.Range("L94").Value = 969.75
$ws.Range("M94").Value = -869.1333
$ws.Range("N94").Value = -1871.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 50499
$ws.Range("J80").Value = 50499
$ws.Range("L80").Value = 50499
$ws.Range("N80").Value = -52745

$ws.Range("H83").Value = 50499
$ws.Range("J83").Value = 50499
$ws.Range("L83").Value = 151497
$ws.Range("N83").Value = -162729

$ws.Range("H122").Value = 1854.6875
$ws.Range("I122").Value = 1882.2174
$ws.Range("K122").Value = 5646.6522
$ws.Range("M122").Value = -3196.6522

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 4278.8335
$ws.Range("I97").Value = 424.5
$ws.Range("J97").Value = 5049.7
$ws.Range("K97").Value = 1273.5
$ws.Range("L97").Value = 15149.1
$ws.Range("M97").Value = -777.5
$ws.Range("N97").Value = -16141.1

$ws.Range("H107").Value = 940.875
$ws.Range("J107").Value = 1370.2142
$ws.Range("L107").Value = 4110.642599999999
$ws.Range("N107").Value = -7950.642599999999

$ws.Range("H113").Value = 63047.5
$ws.Range("I113").Value = 125288.25
$ws.Range("J113").Value = 806.75
$ws.Range("K113").Value = 375864.75
$ws.Range("L113").Value = 2420.25
$ws.Range("M113").Value = -373694.75
$ws.Range("N113").Value = -6760.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 16750.75
$ws.Range("I10").Value = 16750.75
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 16750.75
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -16581.75
$ws.Range("N10").ClearContents()

$ws.Range("H11").Value = 10655.556
$ws.Range("I11").Value = 10557.429
$ws.Range("J11").Value = 10999
$ws.Range("K11").Value = 10557.429
$ws.Range("L11").Value = 10999
$ws.Range("M11").Value = -10418.429
$ws.Range("N11").Value = -11277

$ws.Range("H14").Value = 1437978
$ws.Range("I14").Value = 1667640.1
$ws.Range("J14").Value = 60005
$ws.Range("K14").Value = 1667640.1
$ws.Range("L14").Value = 60005
$ws.Range("M14").Value = -1667472.1
$ws.Range("N14").Value = -60341

$ws.Range("H62").Value = 32000
$ws.Range("I62").Value = 32000
$ws.Range("K62").Value = 32000
$ws.Range("M62").Value = -31314

$ws.Range("H65").Value = 32000
$ws.Range("I65").Value = 32000
$ws.Range("K65").Value = 96000
$ws.Range("M65").Value = -92568

$ws.Range("H122").Value = 3664.7036
$ws.Range("I122").Value = 2283.6191
$ws.Range("K122").Value = 6850.8573
$ws.Range("M122").Value = -4400.8573

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 15000
$ws.Range("I63").Value = 15000
$ws.Range("K63").Value = 15000
$ws.Range("M63").Value = -14251

$ws.Range("H64").Value = 17149
$ws.Range("J64").Value = 17149
$ws.Range("L64").Value = 17149
$ws.Range("N64").Value = -17599

$ws.Range("H66").Value = 15000
$ws.Range("I66").Value = 15000
$ws.Range("K66").Value = 45000
$ws.Range("M66").Value = -41256

$ws.Range("H67").Value = 17149
$ws.Range("J67").Value = 17149
$ws.Range("L67").Value = 17149
$ws.Range("N67").Value = -18709

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 32499.75
$ws.Range("J63").Value = 32499.75
$ws.Range("L63").Value = 32499.75
$ws.Range("N63").Value = -33747.75

$ws.Range("H66").Value = 32499.75
$ws.Range("J66").Value = 32499.75
$ws.Range("L66").Value = 97499.25
$ws.Range("N66").Value = -103739.25

$ws.Range("H126").Value = 3183.1667
$ws.Range("I126").Value = 3239.8
$ws.Range("J126").Value = 2900
$ws.Range("K126").Value = 9719.400000000001
$ws.Range("L126").Value = 8700
$ws.Range("M126").Value = -7249.400000000001
$ws.Range("N126").Value = -13640

$ws.Range("H140").Value = 47500
$ws.Range("J140").Value = 47500
$ws.Range("L140").Value = 47500
$ws.Range("N140").Value = -57860
